$d = $word.ActiveDocument
$c = $d.Content

# wdFindWrap: 1 = wdFindContinue
# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# Replace: 2 = wdReplaceAll

function Replace-All($range, $find, $replace) {
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Header language switcher line -------------------------------------
# The "English" link text and the "live chat" link text are updated via the
# Hyperlinks collection (TextToDisplay) instead of Find/Replace: replacing
# a hyperlink's display run with Find/Replace in this host corrupts the
# run formatting of whatever sits immediately next to the hyperlink, so we
# avoid Find crossing a hyperlink boundary wherever possible.
$h1 = $d.Hyperlinks.Item(1)
$h1.TextToDisplay = "Английский"

# The plain-text run right after the hyperlink (" / Portuguese / ...") is
# replaced leaving its first character (the leading space) untouched, so
# the replacement range no longer starts exactly on the hyperlink boundary
# (which is what triggers the formatting bug in this host).
$p1 = $d.Paragraphs.Item(1)
$afterLink1 = $d.Range($h1.Range.End + 1, $p1.Range.End)
$afterLink1.Text = "Португальский / Французский / Тайский / Вьетнамский / Испанский"

# Second, plain (non-hyperlinked) "English" heading.
Replace-All $c "English" "Английский"

# --- Brief table ---------------------------------------------------------
Replace-All $c "Brief" "Кратко"
Replace-All $c "An email sent to partners in the target country who have sent their documents for review. It will be sent via customer.io" "Электронное письмо, отправленное партнерам в целевой стране, которые отправили свои документы на рассмотрение. Оно будет отправлено через customer.io"
Replace-All $c "Target audience" "Целевая аудитория"
Replace-All $c "Invited partners who have submitted their documents" "Приглашенные партнеры, которые подали свои документы"

# --- Subject line ---------------------------------------------------------
Replace-All $c "Subject line" "Тема"
Replace-All $c " — we got your docs!  " " — мы получили ваши документы!  "

# --- Heading + greeting -----------------------------------------------
Replace-All $c "Thank you for submitting your documents" "Спасибо, что отправили документы"
Replace-All $c "Hi " "Здравствуйте, "

# --- Body paragraphs -----------------------------------------------------
Replace-All $c "Thank you for providing us with your documents for the upcoming " "Благодарим вас за то, что отправили нам документы для предстоящего "
Replace-All $c ". Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation." ". Основываясь на предоставленной вами информации, мы сделаем необходимые приготовления, включая размещение и транспорт."
Replace-All $c "We’re currently reviewing your documents and will reach out to you if we need anything else. " "Сейчас мы изучаем ваши документы и свяжемся с вами, если нам понадобится что-то еще. "

Replace-All $c "If you have any questions, please contact us via " "Если у вас есть вопросы, свяжитесь с нами через "

# "live chat" hyperlink -> "чат"
$h2 = $d.Hyperlinks.Item(2)
$h2.TextToDisplay = "чат"

# " or " between the "live chat"/"WhatsApp" hyperlinks: replace from the
# second character onward so the edit range doesn't begin exactly at the
# hyperlink boundary (same trick as above).
$orProbe = $d.Content.Duplicate
$orProbe.Find.Execute(" or ") | Out-Null
$orTarget = $d.Range($orProbe.Start + 1, $orProbe.End)
$orTarget.Text = "или "

Replace-All $c "If you have any questions, please contact your country manager, " "Если у вас есть вопросы, пожалуйста, свяжитесь с вашим региональным менеджером "
Replace-All $c ", at " " по адресу "
Replace-All $c " or " " или "

Replace-All $c "We look forward to seeing you at " "С нетерпением ждем встречи на "

# --- Comment text ----------------------------------------------------------
# Find.Execute does not operate on comment ranges in this host, so set the
# comment range text directly instead.
$comment = $d.Comments.Item(1)
$comment.Range.Text = "выберите один из вариантов"
